# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 167 of Sheet1 (pushing the
# existing rows 167-195 down to 168-196), and the new row is populated with
# this period's data for "Feria Lagunitas de Puerto Montt" / "Brócoli".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 167, shifting rows 167:195
# down to 168:196 (dimension grows from A1:R195 to A1:R196).
$ws.Rows("167:167").Insert()

# Populate the newly inserted row 167 with the new record.
$ws.Range("A167").Value = 4
$ws.Range("B167").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C167").Value = "Los Lagos"
$ws.Range("D167").Value = 44474
$ws.Range("E167").Value = 10
$ws.Range("F167").Value = 100112023
$ws.Range("G167").Value = "Brócoli"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 1400
$ws.Range("K167").Value = 1100
$ws.Range("L167").Value = 1200
$ws.Range("M167").Value = 1150
$ws.Range("N167").Value = "`$/unidad"
$ws.Range("O167").Value = "Región Metropolitana"
$ws.Range("P167").Value = 1150
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
